$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update ConsentNotGivenReason / ConsentNotGivenDetails for rows 4 and 5 ---
# These cells were blank (number-formatted) and now hold the text "Personal choice",
# styled like the other text cells in the Z/AA columns (e.g. Z2), so first copy the
# formatting from Z2 (style index 7) onto the target cells, then set their values.
$ws.Range("Z2").Copy()
$ws.Range("Z4:AA5").PasteSpecial(-4122)

$ws.Range("Z4").Value2 = "Personal choice"
$ws.Range("AA4").Value2 = "Personal choice"
$ws.Range("Z5").Value2 = "Personal choice"
$ws.Range("AA5").Value2 = "Personal choice"

# --- Update ExpectedFinalMessage wording for rows 4 and 5 ---
$ws.Range("AB4").Value2 = "Consent for the MenACWY vaccination confirmed`nAs you answered ‘yes’ to some of the health questions, we need to check the MenACWY vaccination is suitable for ROSE VOSE. We’ll review your answers and get in touch again soon."
$ws.Range("AB5").Value2 = "Consent for the Td/IPV vaccination confirmed`nSUSAN BYRON is due to get the Td/IPV vaccination at school on"

# --- Update the sheet view's selected cell ---
$ws.Activate()
$ws.Range("AB2").Select()
